$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$ws.Range("K2").Value = "2025-10-28T07:02:58.814490+00:00"
$ws.Range("K3").Value = "2025-10-28T07:03:01.138963+00:00"
$ws.Range("K4").Value = "2025-10-28T07:03:01.138992+00:00"
$ws.Range("K5").Value = "2025-10-28T07:03:01.139011+00:00"
$ws.Range("K6").Value = "2025-10-28T07:03:01.139028+00:00"
$ws.Range("K7").Value = "2025-10-28T07:03:01.139044+00:00"
$ws.Range("K8").Value = "2025-10-28T07:03:03.978942+00:00"
$ws.Range("K9").Value = "2025-10-28T07:03:03.978970+00:00"
$ws.Range("K10").Value = "2025-10-28T07:03:06.720279+00:00"
$ws.Range("K11").Value = "2025-10-28T07:03:09.038582+00:00"
$ws.Range("K12").Value = "2025-10-28T07:03:09.038609+00:00"
$ws.Range("K13").Value = "2025-10-28T07:03:11.373338+00:00"
$ws.Range("K14").Value = "2025-10-28T07:03:11.373366+00:00"
$ws.Range("K15").Value = "2025-10-28T07:03:11.373384+00:00"
$ws.Range("K16").Value = "2025-10-28T07:03:11.373401+00:00"
$ws.Range("K17").Value = "2025-10-28T07:03:19.695236+00:00"
$ws.Range("K18").Value = "2025-10-28T07:03:22.512244+00:00"
$ws.Range("K19").Value = "2025-10-28T07:03:25.538223+00:00"
$ws.Range("K20").Value = "2025-10-28T07:03:28.375744+00:00"
$ws.Range("K21").Value = "2025-10-28T07:03:28.375773+00:00"
$ws.Range("K22").Value = "2025-10-28T07:03:30.712171+00:00"
$ws.Range("K23").Value = "2025-10-28T07:03:30.712199+00:00"
$ws.Range("K24").Value = "2025-10-28T07:03:30.712216+00:00"
$ws.Range("K25").Value = "2025-10-28T07:03:30.712232+00:00"
$ws.Range("K26").Value = "2025-10-28T07:03:30.712248+00:00"
$ws.Range("K27").Value = "2025-10-28T07:03:33.402475+00:00"
$ws.Range("K28").Value = "2025-10-28T07:03:38.413015+00:00"
$ws.Range("K29").Value = "2025-10-28T07:03:38.413042+00:00"
$ws.Range("K30").Value = "2025-10-28T07:03:38.413061+00:00"
$ws.Range("K31").Value = "2025-10-28T07:03:38.413077+00:00"
$ws.Range("K32").Value = "2025-10-28T07:03:40.715210+00:00"
$ws.Range("K33").Value = "2025-10-28T07:03:40.715237+00:00"
$ws.Range("K34").Value = "2025-10-28T07:03:42.996572+00:00"
$ws.Range("K35").Value = "2025-10-28T07:03:42.996600+00:00"
$ws.Range("K36").Value = "2025-10-28T07:03:42.996617+00:00"
$ws.Range("K37").Value = "2025-10-28T07:03:42.996633+00:00"
$ws.Range("K38").Value = "2025-10-28T07:03:42.996649+00:00"
$ws.Range("K39").Value = "2025-10-28T07:03:42.996667+00:00"
$ws.Range("K40").Value = "2025-10-28T07:03:42.996682+00:00"
$ws.Range("K41").Value = "2025-10-28T07:03:42.996697+00:00"
$ws.Range("K42").Value = "2025-10-28T07:03:42.996711+00:00"
$ws.Range("K43").Value = "2025-10-28T07:03:45.358137+00:00"
$ws.Range("K44").Value = "2025-10-28T07:03:45.358166+00:00"
$ws.Range("K45").Value = "2025-10-28T07:03:49.880029+00:00"
$ws.Range("K46").Value = "2025-10-28T07:03:52.123924+00:00"
$ws.Range("K47").Value = "2025-10-28T07:03:52.123955+00:00"
$ws.Range("K48").Value = "2025-10-28T07:03:52.123973+00:00"
$ws.Range("K49").Value = "2025-10-28T07:03:52.123989+00:00"
